$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = "P. point"
$ws.Range("C8").Value = 57
$ws.Range("D8").Value = "'2"
$ws.Range("E8").Value = "Short point (up to 3 mtr.)"
$ws.Range("F8").Value = 256
$ws.Range("G8").Value = "'14592.00"
$ws.Range("I8").Value = "'"

# Row 9
$ws.Range("A9").Value = "P. point"
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = "'3"
$ws.Range("E9").Value = "Medium point (up to 6 mtr.)"
$ws.Range("F9").Value = 472
$ws.Range("G9").Value = "'1888.00"
$ws.Range("I9").Value = "'"

# Row 10
$ws.Range("A10").Value = "Each"
$ws.Range("C10").Value = 24
$ws.Range("D10").Value = "'3.0"
$ws.Range("E10").Value = "P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F10").Value = 23
$ws.Range("G10").Value = "'552.00"
$ws.Range("I10").Value = "'"

# Row 11
$ws.Range("A11").Value = "R. mtr."
$ws.Range("C11").Value = 93
$ws.Range("D11").Value = "'17"
$ws.Range("E11").Value = "25 mm"
$ws.Range("F11").Value = 56
$ws.Range("G11").Value = "'5208.00"
$ws.Range("I11").Value = "'"

# Row 12
$ws.Range("A12").Value = "'"
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = "'15.0"
$ws.Range("E12").Value = "Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = "'0.00"
$ws.Range("I12").Value = "'"

# Row 13
$ws.Range("A13").Value = "Each"
$ws.Range("C13").Value = 66
$ws.Range("D13").Value = "'27"
$ws.Range("E13").Value = "1170mm(+/-10%) LED batten with min. lumen output 2200 lm"
$ws.Range("F13").Value = 492
$ws.Range("G13").Value = "'32472.00"
$ws.Range("I13").Value = "'"

# Row 14
$ws.Range("A14").Value = "'"
$ws.Range("C14").Value = 66
$ws.Range("D14").Value = "'18.0"
$ws.Range("E14").Value = "Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = "'0.00"
$ws.Range("I14").Value = "'"

# Row 15
$ws.Range("A15").Value = "%"
$ws.Range("C15").Value = 17
$ws.Range("D15").Value = "'37"
$ws.Range("E15").Value = "Add Tender Premium "
$ws.Range("G15").Value = "'0.00"
$ws.Range("I15").Value = "'"

# Row 16: remove totals/grand-total row contents, leave only an empty A16
$ws.Range("B16:I16").ClearContents()
$ws.Range("A16").Value = "'"

# Row 17: new Grand Total Rs. row
$ws.Range("A17").Value = "'"
$ws.Range("B17").Value = "'"
$ws.Range("C17").Value = "'"
$ws.Range("D17").Value = "'"
$ws.Range("E17").Value = "Grand Total Rs."
$ws.Range("F17").Value = "'"
$ws.Range("G17").Value = "'54712.00"
$ws.Range("H17").Value = "'54712.00"
$ws.Range("I17").Value = "'"

# Row 18: Tender Premium row (was Grand Total Rs.)
$ws.Range("E18").Value = "Tender Premium @ 0%"
$ws.Range("G18").Value = "'0.00"
$ws.Range("H18").Value = "'0.00"

# Row 19: NET PAYABLE AMOUNT row (was Tender Premium)
$ws.Range("E19").Value = "NET PAYABLE AMOUNT Rs."
$ws.Range("G19").Value = "'54712.00"
$ws.Range("H19").Value = "'54712.00"

# Row 20: delete entirely, sheet now ends at row 19
$ws.Rows(20).Delete()
